# SoIB_summaries.xlsx update
# 1. Rename "High Priority break-up" -> "Interannual update - High Pri"
# 2. Add a new sheet "Major update - High Priority " at the end, containing
#    the data that used to live in "High Priority break-up".
# 3. Update values in "Trends Status", "Priority Status", and
#    "Species qualification" sheets.
# 4. Populate the renamed "Interannual update - High Pri" sheet with its new data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: capture the old "High Priority break-up" data before we overwrite it
# ---------------------------------------------------------------------------
$wsOldBreakup = $wb.Worksheets.Item("High Priority break-up")

# ---------------------------------------------------------------------------
# Step 2: create the new "Major update - High Priority " sheet at the end,
# and copy over the old break-up data (same content, new tab).
# ---------------------------------------------------------------------------
$wsMajor = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsMajor.Name = "Major update - High Priority "

$wsMajor.Range("A1").Value = "Break-up"
$wsMajor.Range("B1").Value = "High Species (no.)"
$wsMajor.Range("C1").Value = "High Species (perc.)"
$wsMajor.Range("D1").Value = "New High Species (no.)"
$wsMajor.Range("E1").Value = "New High Species (perc.)"
$wsMajor.Range("A1:E1").Font.Bold = $true
$wsMajor.Range("A1:E1").HorizontalAlignment = -4108

$wsMajor.Range("A2").Value = "Trend New"
$wsMajor.Range("B2").Value = 4
$wsMajor.Range("C2").Value = 16.7
$wsMajor.Range("D2").Value = 4
$wsMajor.Range("E2").Value = 16.7

$wsMajor.Range("A3").Value = "IUCN"
$wsMajor.Range("B3").Value = 20
$wsMajor.Range("C3").Value = 83.3
$wsMajor.Range("D3").Value = 20
$wsMajor.Range("E3").Value = 83.3

# ---------------------------------------------------------------------------
# Step 3: rename the original sheet and overwrite it with the new
# "Interannual update" data.
# ---------------------------------------------------------------------------
$wsOldBreakup.Name = "Interannual update - High Pri"
$wsInter = $wsOldBreakup

$wsInter.Range("A1").Value = "Break-up"
$wsInter.Range("B1").Value = "High Species (no.)"
$wsInter.Range("C1").Value = "High Species (perc.)"
$wsInter.Range("D1").Value = "New High Species (no.)"
$wsInter.Range("E1").Value = "New High Species (perc.)"

$wsInter.Range("A2").Value = "Trend New"
$wsInter.Range("B2").Value = 63
$wsInter.Range("C2").Value = 61.2
$wsInter.Range("D2").Value = 63
$wsInter.Range("E2").Value = 75

$wsInter.Range("A3").Value = "Trend Different"
$wsInter.Range("B3").Value = 2
$wsInter.Range("C3").Value = 1.9
$wsInter.Range("D3").ClearContents()
$wsInter.Range("E3").ClearContents()

$wsInter.Range("A4").Value = "IUCN"
$wsInter.Range("B4").Value = 38
$wsInter.Range("C4").Value = 36.9
$wsInter.Range("D4").Value = 21
$wsInter.Range("E4").Value = 25

# ---------------------------------------------------------------------------
# Step 4: Trends Status sheet updates
# ---------------------------------------------------------------------------
$wsTrends = $wb.Worksheets.Item("Trends Status")

$wsTrends.Range("B2").Value = 0
$wsTrends.Range("C2").Value = 2
$wsTrends.Range("D2").Value = 0
$wsTrends.Range("E2").Value = 7.7

$wsTrends.Range("B3").Value = 0
$wsTrends.Range("C3").Value = 6
$wsTrends.Range("D3").Value = 0
$wsTrends.Range("E3").Value = 23.1

$wsTrends.Range("B4").Value = 8
$wsTrends.Range("C4").Value = 17
$wsTrends.Range("D4").Value = 44.4
$wsTrends.Range("E4").Value = 65.40000000000001

$wsTrends.Range("B5").Value = 6
$wsTrends.Range("C5").Value = 1
$wsTrends.Range("D5").Value = 33.3
$wsTrends.Range("E5").Value = 3.8

$wsTrends.Range("B6").Value = 4
$wsTrends.Range("C6").Value = 0
$wsTrends.Range("D6").Value = 22.2
$wsTrends.Range("E6").Value = 0

$wsTrends.Range("B7").Value = 46
$wsTrends.Range("C7").Value = 67

$wsTrends.Range("B8").Value = 329
$wsTrends.Range("C8").Value = 300

# ---------------------------------------------------------------------------
# Step 5: Priority Status sheet updates
# ---------------------------------------------------------------------------
$wsPriority = $wb.Worksheets.Item("Priority Status")
$wsPriority.Range("B2").Value = 103
$wsPriority.Range("B3").Value = 286
$wsPriority.Range("B4").Value = 554

# ---------------------------------------------------------------------------
# Step 6: Species qualification sheet updates
# ---------------------------------------------------------------------------
$wsSpecies = $wb.Worksheets.Item("Species qualification")
$wsSpecies.Range("A2").Value = "SoIB Assessment"
$wsSpecies.Range("B2").Value = 393

$wsSpecies.Range("B3").Value = 64
$wsSpecies.Range("C3").Value = 18

$wsSpecies.Range("B4").Value = 93
$wsSpecies.Range("C4").Value = 26

Write-Host "Edit applied"
